# Second commit: add the "DealName/Amount" deal table to Sheet2, then
# leave the selection/active-sheet state the way the author left it
# (Sheet1 at B8, Sheet2 now the active tab with D7 selected).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---- populate Sheet2 with the new deal table ----
$ws2.Range("A1").Value = "DealName"
$ws2.Range("B1").Value = "Amount"

$ws2.Range("A2").Value = "Sun Islam"
$ws2.Range("B2").Value = 100

$ws2.Range("A3").Value = "shobuj Alam200"
$ws2.Range("B3").Value = 200

$ws2.Range("A4").Value = "Badsha Alam"
$ws2.Range("B4").Value = 300

$ws2.Range("A5").Value = "Nazmul "
$ws2.Range("B5").Value = 400

$ws2.Range("A6").Value = "Hasan"
$ws2.Range("B6").Value = 500

$ws2.Range("A7").Value = "Jewal"
$ws2.Range("B7").Value = 600

# column A on Sheet2 was narrowed to fit the new names
$ws2.Columns.Item(1).ColumnWidth = 9.5

# ---- selection / active sheet bookkeeping ----
# Sheet1 is no longer the active tab, but it keeps a remembered selection of B8
$ws1.Range("B8").Select() | Out-Null

# Sheet2 becomes the active tab, with D7 as the selected (empty) cell
$ws2.Activate()
$ws2.Range("D7").Select() | Out-Null
